$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lisez-moi")
$r = $ws.Range("C9")
$insertPos = 14
$newText = "TESTINSERT"
$r.Characters($insertPos, 0).Text = $newText
# Now try to re-apply Bold=false to the newly inserted text AND restore Bold=true on "Source : "
$r.Characters(1,9).Font.Bold = $true
$r.Characters(10, 4+10).Font.Bold = $false
